# Update "江西-漫展信息.xlsx" (展览 / 全部类型 sheets):
#  - Remove the two obsolete events that used to sit in rows 2-3
#    ("九江·ACD动漫游戏嘉年华02" and "江西·樟树静卿国风动漫文化展览会"),
#    shifting every later row up by two.
#  - Renumber the index column (A) back to a clean 1..N sequence.
#  - Refresh the "想去人数" (interest count) for the events whose
#    live counter ticked up since the previous scrape.

$wb = $excel.ActiveWorkbook

# F-column (想去人数) refresh, keyed by the bilibili activity id that is
# unique per row and survives the row shift.
$interestById = @{
    "78362" = 1880  # 赣州·第一届喵喵鱼动漫游戏展
    "80971" = 81    # 万载·第七届馨缘动漫文化展
    "80785" = 194   # 江西·高安首届静卿国风动漫文化展览会
    "81033" = 4479  # 南昌·第一届Cookie动漫嘉年华
    "81207" = 343   # 景德镇·陶溪川×次元文化元宵游园会（ 免费活动）
    "79555" = 1255  # 南昌·meeting动漫游戏嘉年华
    "81362" = 524   # 景德镇·江报国风动漫展
    "81792" = 843   # 江西·ShiningStaR动漫游戏文化节5th
    "81232" = 457   # 南昌·AP动漫游戏嘉年华
    "81691" = 226   # 南昌·CM01动漫游戏博览会
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Drop the two stale leading data rows (old rows 2 & 3); everything
    # below shifts up to close the gap.
    $ws.Range("2:3").EntireRow.Delete()

    $lastRow = $ws.Range("A1").End(4).Row
    for ($r = 2; $r -le $lastRow; $r++) {
        # Clean index numbering: row 2 -> 1, row 3 -> 2, ...
        $ws.Cells.Item($r, 1).Value = $r - 1

        $link = [string]$ws.Cells.Item($r, 8).Value
        foreach ($id in $interestById.Keys) {
            if ($link.EndsWith($id)) {
                $ws.Cells.Item($r, 6).Value = $interestById[$id]
                break
            }
        }
    }
}
